# "Generate Report for Handoff"
# The a748ff53-261c-4aac-aeba-b753486c8eb4.md entry moves from
# "In Translation" to "Ready for handoff" with a refreshed handoff
# timestamp, and its priority flips from "ht" (human translation) to
# "mt" (machine translation), across the Overview summary sheet and the
# per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the a748ff53 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 04:16:04"
$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637

# --- zh-cn sheet: row 3 is the a748ff53 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-13 04:15:54"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.38265482584637

# --- de-de sheet: row 3 is the a748ff53 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-13 04:16:04"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.38265482584637
